# SPS changes in new dashboard test case
# Adds three new mood-based survey columns (Great / OK / Unpleasant) with
# their corresponding SQL lookup queries to the "New Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Dashboard")
$ws.Activate()

# New header cells (row 1)
$ws.Range("O1").Value = "Great"
$ws.Range("P1").Value = "OK"
$ws.Range("Q1").Value = "Unpleasant"

# New query cells (row 2)
$ws.Range("O2").Value = "select count(*) from survey_details where COMPANY_ID=909 and SOURCE not in('Zillow' ,'3rd Party Review') and MOOD like 'Great';"
$ws.Range("P2").Value = "select count(*) from survey_details where COMPANY_ID=909 and SOURCE not in('Zillow' ,'3rd Party Review') and MOOD like 'OK';"
$ws.Range("Q2").Value = "select count(*) from survey_details where COMPANY_ID=909 and SOURCE not in('Zillow' ,'3rd Party Review') and MOOD like 'Unpleasant';"

# Match the narrower custom width used for the new column Q
$ws.Columns.Item(17).ColumnWidth = 8.166666666666666

# Row heights grow to fit the new wrapped query text
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 300

# Scroll the view and select the newly added range, as in the saved file
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("O1:Q2").Select() | Out-Null
